# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-02 (serial 45171) to 2023-09-03 (serial 45172).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C438")
$range.Value = 45172
